$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 141 (Excel shifts existing
# rows 141-158 down to become rows 143-160), matching the new weekly
# data points that were recorded ahead of the previously-first entry.
$ws.Rows.Item(141).Insert()
$ws.Rows.Item(141).Insert()

# New row 141: Maracuyá, Primera, 2023-02-27 (serial 44984)
$ws.Cells.Item(141, 1).Value2 = 1
$ws.Cells.Item(141, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(141, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(141, 4).Value2 = 44984
$ws.Cells.Item(141, 5).Value2 = 15
$ws.Cells.Item(141, 6).Value2 = "Fruta"
$ws.Cells.Item(141, 7).Value2 = 100108
$ws.Cells.Item(141, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(141, 9).Value2 = 100108003
$ws.Cells.Item(141, 10).Value2 = "Maracuyá"
$ws.Cells.Item(141, 11).Value2 = "Sin especificar"
$ws.Cells.Item(141, 12).Value2 = "Primera"
$ws.Cells.Item(141, 13).Value2 = 55
$ws.Cells.Item(141, 14).Value2 = 17000
$ws.Cells.Item(141, 15).Value2 = 18000
$ws.Cells.Item(141, 16).Value2 = 17364
$ws.Cells.Item(141, 17).Value2 = "$/caja 20 kilos"
$ws.Cells.Item(141, 18).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(141, 19).Value2 = 868
$ws.Cells.Item(141, 20).Value2 = 20

# New row 142: Maracuyá, Segunda, 2023-02-27 (serial 44984)
$ws.Cells.Item(142, 1).Value2 = 1
$ws.Cells.Item(142, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(142, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(142, 4).Value2 = 44984
$ws.Cells.Item(142, 5).Value2 = 15
$ws.Cells.Item(142, 6).Value2 = "Fruta"
$ws.Cells.Item(142, 7).Value2 = 100108
$ws.Cells.Item(142, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(142, 9).Value2 = 100108003
$ws.Cells.Item(142, 10).Value2 = "Maracuyá"
$ws.Cells.Item(142, 11).Value2 = "Sin especificar"
$ws.Cells.Item(142, 12).Value2 = "Segunda"
$ws.Cells.Item(142, 13).Value2 = 65
$ws.Cells.Item(142, 14).Value2 = 14000
$ws.Cells.Item(142, 15).Value2 = 15000
$ws.Cells.Item(142, 16).Value2 = 14385
$ws.Cells.Item(142, 17).Value2 = "$/caja 20 kilos"
$ws.Cells.Item(142, 18).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(142, 19).Value2 = 719
$ws.Cells.Item(142, 20).Value2 = 20
